$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 55555756
$ws.Range("I53").Value = 238.27272
$ws.Range("J53").Value = 142857280
$ws.Range("K53").Value = 238.27272
$ws.Range("L53").Value = 142857280
$ws.Range("M53").Value = 398.72728
$ws.Range("N53").Value = -142858554

$ws.Range("H62").Value = 2254.3635
$ws.Range("I62").Value = 1931.125
$ws.Range("J62").Value = 3116.3333
$ws.Range("K62").Value = 1931.125
$ws.Range("L62").Value = 3116.3333
$ws.Range("M62").Value = -1307.125
$ws.Range("N62").Value = -4364.3333

$ws.Range("H64").Value = 3300
$ws.Range("I64").Value = 3300
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3300
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -3052
$ws.Range("N64").ClearContents()

$ws.Range("H65").Value = 2254.3635
$ws.Range("I65").Value = 1931.125
$ws.Range("J65").Value = 3116.3333
$ws.Range("K65").Value = 9655.625
$ws.Range("L65").Value = 15581.6665
$ws.Range("M65").Value = -6535.625
$ws.Range("N65").Value = -21821.6665

$ws.Range("H67").Value = 3300
$ws.Range("I67").Value = 3300
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3300
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -2442
$ws.Range("N67").ClearContents()

$ws.Range("H87").Value = 24082.55
$ws.Range("J87").Value = 24082.55
$ws.Range("L87").Value = 24082.55
$ws.Range("N87").Value = -26578.55

$ws.Range("H90").Value = 24082.55
$ws.Range("J90").Value = 24082.55
$ws.Range("L90").Value = 72247.64999999999
$ws.Range("N90").Value = -84727.64999999999

$ws.Range("H92").Value = 25643020
$ws.Range("I92").Value = 33334928
$ws.Range("K92").Value = 33334928
$ws.Range("M92").Value = -33333680

$ws.Range("H96").Value = 1488.5714
$ws.Range("I96").Value = 1285.6
$ws.Range("J96").Value = 1673.091
$ws.Range("K96").Value = 3856.8
$ws.Range("L96").Value = 5019.272999999999
$ws.Range("M96").Value = -2483.8
$ws.Range("N96").Value = -7765.272999999999

$ws.Range("H129").Value = 1073.0227
$ws.Range("I129").Value = 537.4375
$ws.Range("J129").Value = 1379.0714
$ws.Range("K129").Value = 1612.3125
$ws.Range("L129").Value = 4137.2142
$ws.Range("M129").Value = 3387.6875
$ws.Range("N129").Value = -14137.2142

$ws.Range("H138").Value = 229077.11
$ws.Range("I138").Value = 4273.1113
$ws.Range("J138").Value = 280954.94
$ws.Range("K138").Value = 12819.3339
$ws.Range("L138").Value = 842864.8200000001
$ws.Range("M138").Value = -7679.333899999998
$ws.Range("N138").Value = -853144.8200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3772.7273
$ws.Range("I45").Value = 2900
$ws.Range("K45").Value = 2900
$ws.Range("M45").Value = -2523

$ws.Range("H74").Value = 1591
$ws.Range("I74").Value = 1110.8334
$ws.Range("J74").Value = 2414.1428
$ws.Range("K74").Value = 1110.8334
$ws.Range("L74").Value = 2414.1428
$ws.Range("M74").Value = -236.8334
$ws.Range("N74").Value = -4162.1428

$ws.Range("H77").Value = 1591
$ws.Range("I77").Value = 1110.8334
$ws.Range("J77").Value = 2414.1428
$ws.Range("K77").Value = 5554.166999999999
$ws.Range("L77").Value = 12070.714
$ws.Range("M77").Value = -1186.166999999999
$ws.Range("N77").Value = -20806.714

$ws.Range("H122").Value = 1793.8
$ws.Range("I122").Value = 1886.2
$ws.Range("J122").Value = 1701.4
$ws.Range("K122").Value = 5658.6
$ws.Range("L122").Value = 5104.200000000001
$ws.Range("M122").Value = -3208.6
$ws.Range("N122").Value = -10004.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4301.387
$ws.Range("I31").Value = 825.9375
$ws.Range("J31").Value = 8008.533
$ws.Range("K31").Value = 825.9375
$ws.Range("L31").Value = 8008.533
$ws.Range("M31").Value = -530.9375
$ws.Range("N31").Value = -8598.532999999999

$ws.Range("H33").Value = 3000
$ws.Range("I33").Value = 3000
$ws.Range("K33").Value = 3000
$ws.Range("M33").Value = -2621

$ws.Range("H34").Value = 4301.387
$ws.Range("I34").Value = 825.9375
$ws.Range("J34").Value = 8008.533
$ws.Range("K34").Value = 825.9375
$ws.Range("L34").Value = 8008.533
$ws.Range("M34").Value = -623.9375
$ws.Range("N34").Value = -8412.532999999999

$ws.Range("H36").Value = 4988
$ws.Range("I36").Value = 4988
$ws.Range("K36").Value = 4988
$ws.Range("M36").Value = -4600

$ws.Range("H40").Value = 4988
$ws.Range("I40").Value = 4988
$ws.Range("K40").Value = 4988
$ws.Range("M40").Value = -4828

$ws.Range("H58").Value = 1517.0714
$ws.Range("I58").Value = 1258.75
$ws.Range("J58").Value = 1620.4
$ws.Range("K58").Value = 1258.75
$ws.Range("L58").Value = 1620.4
$ws.Range("M58").Value = -1055.75
$ws.Range("N58").Value = -2026.4

$ws.Range("H63").Value = 100000
$ws.Range("J63").Value = 100000
$ws.Range("L63").Value = 100000
$ws.Range("N63").Value = -101372

$ws.Range("H66").Value = 100000
$ws.Range("J66").Value = 100000
$ws.Range("L66").Value = 300000
$ws.Range("N66").Value = -306864

$ws.Range("H122").Value = 1422.6809
$ws.Range("I122").Value = 1019.0476
$ws.Range("J122").Value = 1748.6923
$ws.Range("K122").Value = 3057.1428
$ws.Range("L122").Value = 5246.0769
$ws.Range("M122").Value = -607.1428000000001
$ws.Range("N122").Value = -10146.0769

$ws.Range("H134").Value = 3174.5
$ws.Range("I134").Value = 1437.3334
$ws.Range("J134").Value = 4477.375
$ws.Range("K134").Value = 4312.0002
$ws.Range("L134").Value = 13432.125
$ws.Range("M134").Value = -1777.0002
$ws.Range("N134").Value = -18502.125

$ws.Range("H136").Value = 1517.0714
$ws.Range("I136").Value = 1258.75
$ws.Range("J136").Value = 1620.4
$ws.Range("K136").Value = 3776.25
$ws.Range("L136").Value = 4861.200000000001
$ws.Range("M136").Value = -1226.25
$ws.Range("N136").Value = -9961.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 45015.145
$ws.Range("I21").Value = 506
$ws.Range("J21").Value = 52433.332
$ws.Range("K21").Value = 506
$ws.Range("L21").Value = 52433.332
$ws.Range("M21").Value = -333
$ws.Range("N21").Value = -52779.332

$ws.Range("H30").Value = 45015.145
$ws.Range("I30").Value = 506
$ws.Range("J30").Value = 52433.332
$ws.Range("K30").Value = 506
$ws.Range("L30").Value = 52433.332
$ws.Range("M30").Value = -401
$ws.Range("N30").Value = -52643.332

$ws.Range("H122").Value = 3934.1667
$ws.Range("I122").Value = 2616.524
$ws.Range("K122").Value = 7849.572
$ws.Range("M122").Value = -5399.572

$ws.Range("H132").Value = 2233.3809
$ws.Range("I132").Value = 1556.3572
$ws.Range("J132").Value = 3587.4285
$ws.Range("K132").Value = 4669.071599999999
$ws.Range("L132").Value = 10762.2855
$ws.Range("M132").Value = -2139.071599999999
$ws.Range("N132").Value = -15822.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 8012.75
$ws.Range("I33").Value = 8000
$ws.Range("J33").Value = 8017
$ws.Range("K33").Value = 8000
$ws.Range("L33").Value = 8017
$ws.Range("M33").Value = -7710
$ws.Range("N33").Value = -8597

$ws.Range("H68").Value = 2000
$ws.Range("J68").Value = 2000
$ws.Range("L68").Value = 2000
$ws.Range("N68").Value = -3498

$ws.Range("H71").Value = 2000
$ws.Range("J71").Value = 2000
$ws.Range("L71").Value = 10000
$ws.Range("N71").Value = -17488

$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2392.8235
$ws.Range("I122").Value = 2261.6365
$ws.Range("J122").Value = 2633.3333
$ws.Range("K122").Value = 6784.9095
$ws.Range("L122").Value = 7899.999899999999
$ws.Range("M122").Value = -4334.9095
$ws.Range("N122").Value = -12799.9999

$ws.Range("H136").Value = 4976.2666
$ws.Range("I136").Value = 5014.9
$ws.Range("J136").Value = 5014.9
$ws.Range("K136").Value = 15044.7
$ws.Range("L136").Value = 14697
$ws.Range("M136").Value = -12494.7
$ws.Range("N136").Value = -19797
